$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fix the misspelling "mecahnics" -> "mechanics" (Game Mechanics
#    intro paragraph).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "mecahnics", $false, $false, $false, $false, $false,
    $true, 1, $false, "mechanics", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Grammar fix "a artificial intelligence" -> "an artificial
#    intelligence" (AI Engine intro paragraph).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "to create a artificial intelligence", $false, $false, $false, $false, $false,
    $true, 1, $false, "to create an artificial intelligence", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Rewrite the "This GUI will be created..." sentence.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "This GUI will be created using a c++ graphics library, possibly FLTK.",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "This GUI will be created using QT, a c++ library.", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Server section rework.
#    The paragraph that used to describe HandleMove() now describes
#    HandleConnection(), and the four paragraphs that followed it
#    (HandleUndo/Redo, HandleDifficulty, HandleExit, HandleDisplay)
#    are removed entirely.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "HandleMove()- the function used if the user has input a move, with call on game mechanics functions.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "HandleConnection()- this function is the actual flow for the gameplay, it handles all the client/server interactions. The design is a little cluttered, it could have been broken up into more functions.",
    2) | Out-Null

# Locate that paragraph again (now renamed) so we can relocate the
# hidden "_GoBack" bookmark to sit right after it, and delete the
# four obsolete paragraphs that follow.
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "HandleConnection*") {
        $handleConnectionIndex = $i
        break
    }
}

# Move the "_GoBack" bookmark from the start of "The server is the
# central point..." paragraph to the end of the HandleConnection
# paragraph (right before its paragraph mark).
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

$hcPara = $d.Paragraphs.Item($handleConnectionIndex)
$hcRange = $hcPara.Range
$bmPos = $hcRange.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Delete the four now-obsolete paragraphs right after the
# HandleConnection paragraph.
for ($n = 0; $n -lt 4; $n++) {
    $p = $d.Paragraphs.Item($handleConnectionIndex + 1)
    $p.Range.Delete() | Out-Null
}
